# Update "想去人数" (want-to-go count) values in F column across sheets,
# reflecting refreshed data pulled from the source site.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 400
$ws.Range("F5").Value = 947
$ws.Range("F6").Value = 5038
$ws.Range("F7").Value = 402
$ws.Range("F8").Value = 581
$ws.Range("F9").Value = 880
$ws.Range("F16").Value = 1655
$ws.Range("F18").Value = 739
$ws.Range("F20").Value = 178
$ws.Range("F22").Value = 476
$ws.Range("F24").Value = 1040
$ws.Range("F27").Value = 2035
$ws.Range("F29").Value = 84
$ws.Range("F31").Value = 213
$ws.Range("F37").Value = 573
$ws.Range("F38").Value = 75
$ws.Range("F40").Value = 30

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 137

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 400
$ws.Range("F6").Value = 947
$ws.Range("F8").Value = 5038
$ws.Range("F9").Value = 402
$ws.Range("F10").Value = 581
$ws.Range("F12").Value = 137
$ws.Range("F13").Value = 880
$ws.Range("F23").Value = 1655
$ws.Range("F25").Value = 739
$ws.Range("F27").Value = 178
$ws.Range("F30").Value = 476
$ws.Range("F32").Value = 1040
$ws.Range("F34").Value = 2035
$ws.Range("F36").Value = 84
$ws.Range("F38").Value = 213
$ws.Range("F43").Value = 573
$ws.Range("F44").Value = 75
$ws.Range("F46").Value = 30
